$d = $word.ActiveDocument

$replacements = @(
    @{old="343×6="; new="376×7="},
    @{old="760×8="; new="787×8="},
    @{old="633×5="; new="836×3="},
    @{old="938×3="; new="573×3="},
    @{old="367×4="; new="106×8="},
    @{old="224×9="; new="326×6="},
    @{old="475×7="; new="592×9="},
    @{old="664×5="; new="757×4="},
    @{old="267×3="; new="609×5="},
    @{old="864×4="; new="901×7="},
    @{old="441×7="; new="188×3="},
    @{old="519×5="; new="740×4="},
    @{old="460×8="; new="587×7="},
    @{old="826×7="; new="579×9="},
    @{old="668×4="; new="475×5="},
    @{old="128×8="; new="159×5="},
    @{old="848×2="; new="376×6="},
    @{old="332×5="; new="743×4="},
    @{old="255×4="; new="847×3="},
    @{old="530×7="; new="881×2="},
    @{old="879×5="; new="224×7="},
    @{old="597×8="; new="698×8="},
    @{old="499×7="; new="161×5="},
    @{old="429×9="; new="580×7="},
    @{old="831×7="; new="151×4="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
